$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for rows 3 and 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 04:24:47"
$wsZhCn.Range("H3").Value = "2016-03-23 04:25:13"
$wsZhCn.Range("E5").Value = "2016-03-23 04:24:47"
$wsZhCn.Range("H5").Value = "2016-03-23 04:25:13"

# de-de sheet: update Correspond Handoff/Handback Datetime for rows 3 and 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 04:24:51"
$wsDeDe.Range("H3").Value = "2016-03-23 04:25:20"
$wsDeDe.Range("E5").Value = "2016-03-23 04:24:51"
$wsDeDe.Range("H5").Value = "2016-03-23 04:25:20"
